$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 96
$ws.Range("I2").Value = 244
$ws.Range("J2").Value = 1097
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 295
$ws.Range("M2").Value = 23
$ws.Range("N2").Value = 184
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 2
$ws.Range("S2").Value = 133
$ws.Range("T2").Value = 176
$ws.Range("U2").Value = 10
$ws.Range("V2").Value = 1629
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1590
$ws.Range("Z2").Value = 17
$ws.Range("AA2").Value = 9
